# Update the "cutoff freq" input on the "Branch LPF" worksheet from 1800 Hz
# to 1400 Hz. All dependent formulas on that sheet (radian cutoff, prewarp,
# a1/b0/b1 coefficients, scaled coefficients, etc.) recalculate automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Branch LPF")
$ws.Range("B2").Value = 1400
